$wb = $excel.ActiveWorkbook

$wsTotaal = $wb.Worksheets.Item("totaal")
$wsWeek47 = $wb.Worksheets.Item("weeknr 47")

# --- Fill in the new log entry on row 8 of "weeknr 47" ---
$wsWeek47.Range("B8").Value = 41235
$wsWeek47.Range("C8").Value = 0.36458333333333331
$wsWeek47.Range("D8").Value = 0.5
$wsWeek47.Range("E8").Value = 2
$wsWeek47.Range("F8").Value = "activatie proces anngemaakt, wachtword word gewijzigd bij administratie"
$wsWeek47.Range("G8").Formula = "=D8-C8"

# Clear out the stale shared-formula results that used to fill A9:A23,
# while re-establishing the shared formula anchored at A8 (spanning
# A8:A23) before clearing the follower cells again.
$wsWeek47.Range("A9:A23").ClearContents()
$wsWeek47.Range("A8:A23").Formula = "=B8"
$wsWeek47.Range("A9:A23").ClearContents()

# --- Update selections / active sheet/tab ---
# Make "weeknr 47" selection land on F17, then switch to "totaal" and
# select C16 there, so "totaal" ends up the active (last-selected) sheet.
$wsWeek47.Range("F17").Select()
$wsTotaal.Range("C16").Select()
